{"js": "// Turn off \"page break before\" on every body paragraph and on every\n// built-in heading/title style (Heading 1-6, Title, Subtitle), making the\n// w:pageBreakBefore setting explicit (val=\"0\") instead of inherited/unset.\n\n// 1) Body paragraphs: set paragraphFormat.pageBreakBefore = false on each.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nfor (const paragraph of paragraphs.items) {\n  paragraph.paragraphFormat.pageBreakBefore = false;\n}\nawait context.sync();\n\n// 2) Paragraph styles that carry keepNext/keepLines (headings + title/subtitle).\nconst styleNames = [\n  \"Heading 1\",\n  \"Heading 2\",\n  \"Heading 3\",\n  \"Heading 4\",\n  \"Heading 5\",\n  \"Heading 6\",\n  \"Title\",\n  \"Subtitle\",\n];\n\nfor (const name of styleNames) {\n  const style = context.document.getStyles().getByNameOrNullObject(name);\n  style.paragraphFormat.pageBreakBefore = false;\n}\nawait context.sync();\n", "ps1": "# Turn off \"page break before\" on every body paragraph and on every\n# built-in heading/title style (Heading 1-6, Title, Subtitle), making the\n# PageBreakBefore setting explicit (False / w:val=\"0\") instead of\n# inherited/unset.\n\n$d = $word.ActiveDocument\n\n# 1) Body paragraphs: set PageBreakBefore = False on each paragraph.\nforeach ($p in $d.Paragraphs) {\n    $p.PageBreakBefore = 0\n}\n\n# 2) Paragraph styles that carry KeepNext/KeepLines (headings + title/subtitle).\n$styleNames = @(\"Heading1\", \"Heading2\", \"Heading3\", \"Heading4\", \"Heading5\", \"Heading6\", \"Title\", \"Subtitle\")\nforeach ($name in $styleNames) {\n    $s = $d.Styles($name)\n    $s.ParagraphFormat.PageBreakBefore = 0\n}\n"}
